$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - MAZA RIOFRIO CINTHIA NATELAHI stays, value updated
$ws.Range("B2").Value = 124

# Row 3 - now PALACIOS PANTA LUIS MIGUEL (was PANTA NIMA FREDDY ROLAND JUNIOR)
$ws.Range("A3").Value = "PALACIOS PANTA LUIS MIGUEL"
$ws.Range("B3").Value = 120

# Row 4 - now PANTA NIMA FREDDY ROLAND JUNIOR (was PALACIOS PANTA LUIS MIGUEL)
$ws.Range("A4").Value = "PANTA NIMA FREDDY ROLAND JUNIOR"
$ws.Range("B4").Value = 118

# Row 5 - VEGA ZAPATA JESUS GABRIEL stays, value updated
$ws.Range("B5").Value = 117

# Row 6 - now PANTA VARONA CANDY ELIZABETH (was ELIAS MACHADO JUANA MARGOT)
$ws.Range("A6").Value = "PANTA VARONA CANDY ELIZABETH"
$ws.Range("B6").Value = 116

# Row 7 - CRISANTO CARMEN ROSITA ABIGAIL stays, value updated
$ws.Range("B7").Value = 115

# Row 8 - now ELIAS MACHADO JUANA MARGOT (was PANTA VARONA CANDY ELIZABETH)
$ws.Range("A8").Value = "ELIAS MACHADO JUANA MARGOT"
$ws.Range("B8").Value = 115

# Row 9 - SALAZAR VEGA MARIA FERNANDA stays, value updated
$ws.Range("B9").Value = 114

# Row 10 - HIDALGO MOSCOL YESSICA JAZMIN stays, value updated
$ws.Range("B10").Value = 107

# Row 11 - TALLEDO ELIAS ANDREA ALESSANDRA stays, value updated
$ws.Range("B11").Value = 102
